# Add data for 2021-12-10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab title reflects the new "through" date)
$ws.Name = "Through 2021-12-02"

# Update the December row label
$ws.Range("A14").Value = "December (through 12-02)"

# Row 14 - December monthly figures
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 0.1667
$ws.Range("I14").Value = 6
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 6
$ws.Range("M14").Value = 0.1429
$ws.Range("O14").Value = 2
$ws.Range("R14").Value = 8
$ws.Range("U14").Value = 14

# Row 15 - Totals
$ws.Range("E15").Value = 61
$ws.Range("F15").Value = 508
$ws.Range("G15").Value = 0.1072
$ws.Range("I15").Value = 764
$ws.Range("J15").Value = 0.0762
$ws.Range("K15").Value = 75
$ws.Range("L15").Value = 614
$ws.Range("M15").Value = 0.1089
$ws.Range("O15").Value = 482
$ws.Range("P15").Value = 0.1007
$ws.Range("R15").Value = 1208
$ws.Range("S15").Value = 0.0503
$ws.Range("U15").Value = 1561
$ws.Range("V15").Value = 0.0585
